# Re-shuffle the tied-ranking category labels produced by the (re-run)
# serial analysis generation. Only the category-name column (A) changes;
# the numeric values in column B stay exactly as they were since these
# rows are ties in the ranking.

$wb = $excel.ActiveWorkbook

# --- sheet "max-arrecad" ---------------------------------------------
$wsMax = $wb.Worksheets.Item("max-arrecad")

$wsMax.Range("A3").Value  = "midia_independente"
$wsMax.Range("A4").Value  = "disputa"
$wsMax.Range("A5").Value  = "herois"
$wsMax.Range("A6").Value  = "terror"
$wsMax.Range("A7").Value  = "politica"
$wsMax.Range("A8").Value  = "religiosidade"
$wsMax.Range("A9").Value  = "humor"
$wsMax.Range("A10").Value = "erotismo"

$wsMax.Range("A16").Value = "angelo_agostini"
$wsMax.Range("A17").Value = "ccxp"
$wsMax.Range("A19").Value = "hqmix"

# --- sheet "tx-sucesso" -----------------------------------------------
$wsTx = $wb.Worksheets.Item("tx-sucesso")

$wsTx.Range("A9").Value  = "lgbtqiamais"
$wsTx.Range("A10").Value = "zine"
